$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2023-08-01 Tuesday" "2023-08-02 Wednesday"

Replace-Text "35÷6=5, 5" "36÷6=6, 0"
Replace-Text "23÷4=5, 3" "82÷4=20, 2"
Replace-Text "67÷3=22, 1" "92÷8=11, 4"
Replace-Text "16÷2=8, 0" "67÷9=7, 4"
Replace-Text "50÷3=16, 2" "56÷5=11, 1"
Replace-Text "32÷3=10, 2" "48÷5=9, 3"
Replace-Text "68÷7=9, 5" "52÷6=8, 4"
Replace-Text "75÷5=15, 0" "50÷4=12, 2"
Replace-Text "10÷5=2, 0" "93÷3=31, 0"
Replace-Text "44÷9=4, 8" "66÷8=8, 2"
Replace-Text "75÷7=10, 5" "98÷5=19, 3"
Replace-Text "53÷3=17, 2" "97÷6=16, 1"
Replace-Text "92÷5=18, 2" "37÷8=4, 5"
Replace-Text "91÷5=18, 1" "52÷5=10, 2"
Replace-Text "71÷9=7, 8" "76÷9=8, 4"
Replace-Text "56÷3=18, 2" "67÷5=13, 2"
Replace-Text "72÷3=24, 0" "44÷5=8, 4"
Replace-Text "81÷3=27, 0" "61÷8=7, 5"
Replace-Text "15÷4=3, 3" "15÷3=5, 0"
Replace-Text "33÷8=4, 1" "87÷9=9, 6"
Replace-Text "44÷4=11, 0" "31÷2=15, 1"
Replace-Text "23÷7=3, 2" "24÷2=12, 0"
Replace-Text "32÷7=4, 4" "65÷4=16, 1"
Replace-Text "55÷7=7, 6" "86÷4=21, 2"
Replace-Text "90÷3=30, 0" "47÷5=9, 2"
